$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that currently sits right
#    after the title run ("PI SIGMA ALPHA, Alpha Pi Chapter").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
    Write-Host "Removed old _GoBack bookmark"
}

# ------------------------------------------------------------------
# 2. Bump the lifetime membership fee from $45.00 to $35.00.
# ------------------------------------------------------------------
$okFee = $d.Content.Find.Execute(
    "a lifetime membership fee of `$45.00",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a lifetime membership fee of `$35.00", 2)
Write-Host "Fee replace ok:" $okFee

# ------------------------------------------------------------------
# 3. Update the drop-off instructions: replace the "main office ..."
#    sentence fragment with the faculty sponsor's info.
# ------------------------------------------------------------------
$okAddr = $d.Content.Find.Execute(
    "PI SIGMA ALPHA to the main office of the Political Science Department in 104 Baldwin Hall and leave it with receptionist. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "PI SIGMA ALPHA to the faculty sponsor, Brittany Bramlett, in 380C Baldwin Hall. ", 2)
Write-Host "Address replace ok:" $okAddr

# ------------------------------------------------------------------
# 4. Re-add the "_GoBack" bookmark at the very end of that same
#    paragraph (right after "...380C Baldwin Hall. ").
#    A collapsed range positioned exactly at the last character of a
#    paragraph can mis-resolve in this host, so we temporarily append
#    a one-character marker, wrap the bookmark around it, then clear
#    the marker text again - leaving a correctly placed empty
#    bookmark immediately before the closing paragraph mark.
# ------------------------------------------------------------------
$r = $d.Content
$okFind = $r.Find.Execute("380C Baldwin Hall. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Locate end of sentence ok:" $okFind
$r.Collapse(0)
$r.InsertAfter("@")
$bm = $d.Range($r.Start, $r.End)
$d.Bookmarks.Add("_GoBack", $bm)
$bm.Text = ""

Write-Host "Edits applied"
